$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to reflect the login test data it now holds
$ws.Name = "ValidLogin"

# Write the data row first (valid login credentials used to drive the test),
# then the header row - matches the order the original author typed them in,
# which determines the shared-string table ordering.
$ws.Range("A2").Value = "admin"
$ws.Range("B2").Value = "pointofsale"
$ws.Range("B1").Value = "Password"
$ws.Range("A1").Value = "Username"

# Autofit columns so their widths match content (bestFit) like Excel would do
$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(2).AutoFit() | Out-Null
